$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet (tab + workbook.xml <sheet name=...>)
$ws.Name = "GossA"

# 2. Tiny precision corrections to a handful of row-13 values
$ws.Cells.Item(13, 4).Value  = 0.9933137835985731   # D13
$ws.Cells.Item(13, 8).Value  = 0.9933137835985731   # H13
$ws.Cells.Item(13, 12).Value = 0.9920335538727214   # L13
$ws.Cells.Item(13, 14).Value = 0.9943668709229161   # N13

# 3. Append a new row 16 ("HexGrid-60degTilt5degRes" / index 14), matching
#    the look (border + bold + centered) of the other index/label rows.
$ws.Range("A15").Copy()
$ws.Range("A16").PasteSpecial(-4122)

$ws.Cells.Item(16, 1).Value  = 14
$ws.Cells.Item(16, 2).Value  = "HexGrid-60degTilt5degRes"
$ws.Cells.Item(16, 3).Value  = 1.201031381252061
$ws.Cells.Item(16, 4).Value  = 0.9924703951312377
$ws.Cells.Item(16, 5).Value  = 0.9834174039073463
$ws.Cells.Item(16, 6).Value  = 0.9305288542310139
$ws.Cells.Item(16, 7).Value  = 1.201031381252061
$ws.Cells.Item(16, 8).Value  = 0.9924703951312377
$ws.Cells.Item(16, 9).Value  = 1.036209067011093
$ws.Cells.Item(16, 10).Value = 0.9172375850247387
$ws.Cells.Item(16, 11).Value = 1.026900555561042
$ws.Cells.Item(16, 12).Value = 0.9346251806368507
$ws.Cells.Item(16, 13).Value = 1.201031381252061
$ws.Cells.Item(16, 14).Value = 0.987943899519292
$ws.Cells.Item(16, 15).Value = 1.026862008630415
$ws.Cells.Item(16, 16).Value = 1.002802552844423
